# Actualización al 08 de octubre 2023
# Adds new contribution / expense rows and a brand-new "Actividad GOAT" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Ingreso" sheet — 3 new rows (aportes del 08/10/2023)
# ---------------------------------------------------------------------------
$ingreso = $wb.Worksheets.Item("Ingreso")

$ingreso.Range("A542").Value = 45207
$ingreso.Range("B542").Value = "Wilkin"
$ingreso.Range("C542").Value = 100
$ingreso.Range("D542").Value = "Aporte"

$ingreso.Range("A543").Value = 45207
$ingreso.Range("B543").Value = "Fernando"
$ingreso.Range("C543").Value = 100
$ingreso.Range("D543").Value = "Aporte"

$ingreso.Range("A544").Value = 45207
$ingreso.Range("B544").Value = "Javier"
$ingreso.Range("C544").Value = -60
$ingreso.Range("D544").Value = "Préstamo"

# ---------------------------------------------------------------------------
# 2. "Gastos" sheet — 3 new rows (gastos del 08/10/2023)
# ---------------------------------------------------------------------------
$gastos = $wb.Worksheets.Item("Gastos")

$gastos.Range("A63").Value = 45207
$gastos.Range("B63").Value = "Arbitro, agua y hielo"
$gastos.Range("C63").Formula = "=800+260"

$gastos.Range("A64").Value = 45207
$gastos.Range("B64").Value = "Adelanto Arbitro"
$gastos.Range("C64").Value = 400

$gastos.Range("A65").Value = 45207
$gastos.Range("B65").Value = "Junte GOAT"
$gastos.Range("C65").Value = 9833

# ---------------------------------------------------------------------------
# 3. New sheet "Actividad GOAT" at the end of the workbook
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$goat = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$goat.Name = "Actividad GOAT"

$goat.Range("A1").Value = "Carbón"
$goat.Range("B1").Value = 220

$goat.Range("A2").Value = "Colmado"
$goat.Range("B2").Value = 230

$goat.Range("A3").Value = "Compra"
$goat.Range("B3").Formula = "=1668+1200+3400+1625"

$goat.Range("A4").Value = "Empanadas"
$goat.Range("B4").Value = 840

$goat.Activate() | Out-Null
$goat.Range("B1:B4").Select() | Out-Null
$goat.Range("B4").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Restore selections / active sheet to match the saved workbook state
# ---------------------------------------------------------------------------
$gastos.Activate() | Out-Null
$gastos.Range("D68").Select() | Out-Null

$ingreso.Activate() | Out-Null
$ingreso.Range("B549").Select() | Out-Null

Write-Output "done"
